# restructure for DBServer fixed bugs for Lua's interfaces
#
# The "AIServer" row (row 16) in the Server sheet is renamed to "DBServer":
#   A16: "AIServer_1" -> "DBServer_1"
#   C16: "AIServer"   -> "DBServer"
# Leave every other cell/value/format untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = "DBServer_1"
$ws.Range("C16").Value = "DBServer"

# Author ended up with the selection resting on the renamed cell (C16).
$ws.Range("C16").Select()
